# This sheet is a weekly price log (Hortaliza / Acelga, Macroferia Regional
# de Talca). The commit adds one new week's observation at the top of the
# data block (row 39) and pushes every later weekly record down by one row,
# through row 165; the record that used to occupy row 165 is preserved by
# appending it as a brand-new row 166 at the end of the table.
#
# Only the "weekly" columns move: D (Fecha), J (Volumen), K/L/M (Precio
# minimo/maximo/promedio ponderado), O (Origen) and P (Precio $/Kg). The
# constant descriptive columns (A,B,C,E,F,G,H,I,N,Q,R) never change for this
# subset, so row 166 is simplest to build as a full copy of (old) row 165.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 39
$lastDataRow  = 165
$newLastRow   = 166

# --- 1. Append the new row 166 = a full copy of the (still unmodified) row 165 ---
$lastRowValues = $ws.Range("A$lastDataRow`:R$lastDataRow").Value2
$ws.Range("A$newLastRow`:R$newLastRow").Value2 = $lastRowValues
# Value2 does not carry number formats across, so the date column needs it
# copied explicitly to keep the "YYYY-MM-DD HH:MM:SS" display.
$ws.Range("D$newLastRow").NumberFormat = $ws.Range("D$lastDataRow").NumberFormat

# --- 2. Snapshot the "before" state of the moving columns for rows 39..165 ---
$old = $ws.Range("A$firstDataRow`:R$lastDataRow").Value2

# --- 3. Shift rows 40..165 down: row r takes what row (r-1) used to hold ---
for ($r = $lastDataRow; $r -ge ($firstDataRow + 1); $r--) {
    $srcIdx = ($r - 1) - $firstDataRow + 1   # 1-based row index into $old

    $ws.Cells.Item($r, 4).Value2  = $old[$srcIdx, 4]   # D - Fecha
    $ws.Cells.Item($r, 10).Value2 = $old[$srcIdx, 10]  # J - Volumen
    $ws.Cells.Item($r, 11).Value2 = $old[$srcIdx, 11]  # K - Precio minimo
    $ws.Cells.Item($r, 12).Value2 = $old[$srcIdx, 12]  # L - Precio maximo
    $ws.Cells.Item($r, 13).Value2 = $old[$srcIdx, 13]  # M - Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value2 = $old[$srcIdx, 15]  # O - Origen
    $ws.Cells.Item($r, 16).Value2 = $old[$srcIdx, 16]  # P - Precio $/Kg
}

# --- 4. Row 39 becomes the newly reported week ---
$ws.Cells.Item($firstDataRow, 4).Value2  = 44453
$ws.Cells.Item($firstDataRow, 10).Value2 = 400
$ws.Cells.Item($firstDataRow, 11).Value2 = 2500
$ws.Cells.Item($firstDataRow, 12).Value2 = 2500
$ws.Cells.Item($firstDataRow, 13).Value2 = 2500
$ws.Cells.Item($firstDataRow, 16).Value2 = 625
